$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 418.94736
$ws.Range("I19").Value = 415.54544
$ws.Range("K19").Value = 415.54544
$ws.Range("M19").Value = -240.54544
$ws.Range("H76").Value = 5698.875
$ws.Range("J76").Value = 6133.5
$ws.Range("L76").Value = 6133.5
$ws.Range("N76").Value = -6763.5
$ws.Range("H79").Value = 5698.875
$ws.Range("J79").Value = 6133.5
$ws.Range("L79").Value = 6133.5
$ws.Range("N79").Value = -8317.5
$ws.Range("H111").Value = 2815.889
$ws.Range("I111").Value = 2643.1667
$ws.Range("J111").Value = 3161.3333
$ws.Range("K111").Value = 7929.500100000001
$ws.Range("L111").Value = 9483.999899999999
$ws.Range("M111").Value = -4862.500100000001
$ws.Range("N111").Value = -15617.9999
$ws.Range("H120").Value = 77000
$ws.Range("J120").Value = 77000
$ws.Range("L120").Value = 77000
$ws.Range("N120").Value = -86676
$ws.Range("H132").Value = 31252902
$ws.Range("I132").Value = 34485510
$ws.Range("K132").Value = 103456530
$ws.Range("M132").Value = -103454000
$ws.Range("H135").Value = 3089.4348
$ws.Range("I135").Value = 2865.5715
$ws.Range("J135").Value = 5440
$ws.Range("K135").Value = 25790.1435
$ws.Range("L135").Value = 48960
$ws.Range("M135").Value = -23255.1435
$ws.Range("N135").Value = -54030
$ws.Range("H138").Value = 3904.8044
$ws.Range("I138").Value = 1894.15
$ws.Range("J138").Value = 5451.4614
$ws.Range("K138").Value = 5682.450000000001
$ws.Range("L138").Value = 16354.3842
$ws.Range("M138").Value = -542.4500000000007
$ws.Range("N138").Value = -26634.3842
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2964.0815
$ws.Range("I32").Value = 2387.4285
$ws.Range("K32").Value = 2387.4285
$ws.Range("M32").Value = -2100.4285
$ws.Range("H132").Value = 18757
$ws.Range("I132").Value = 11382.95
$ws.Range("K132").Value = 34148.85000000001
$ws.Range("M132").Value = -31618.85000000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3352.96
$ws.Range("I86").Value = 3089.6365
$ws.Range("J86").Value = 3559.8572
$ws.Range("K86").Value = 3089.6365
$ws.Range("L86").Value = 3559.8572
$ws.Range("M86").Value = -1966.6365
$ws.Range("N86").Value = -5805.8572
$ws.Range("H89").Value = 3352.96
$ws.Range("I89").Value = 3089.6365
$ws.Range("J89").Value = 3559.8572
$ws.Range("K89").Value = 15448.1825
$ws.Range("L89").Value = 17799.286
$ws.Range("M89").Value = -9832.182500000001
$ws.Range("N89").Value = -29031.286
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2661.9375
$ws.Range("J31").Value = 3188
$ws.Range("L31").Value = 3188
$ws.Range("N31").Value = -3778
$ws.Range("H34").Value = 2661.9375
$ws.Range("J34").Value = 3188
$ws.Range("L34").Value = 3188
$ws.Range("N34").Value = -3592
$ws.Range("H58").Value = 2630
$ws.Range("I58").Value = 2511.05
$ws.Range("K58").Value = 2511.05
$ws.Range("M58").Value = -2308.05
$ws.Range("H99").Value = 4533.9287
$ws.Range("I99").Value = 3697.6667
$ws.Range("J99").Value = 5161.125
$ws.Range("K99").Value = 3697.6667
$ws.Range("L99").Value = 5161.125
$ws.Range("M99").Value = -2199.6667
$ws.Range("N99").Value = -8157.125
$ws.Range("H126").Value = 4533.9287
$ws.Range("I126").Value = 3697.6667
$ws.Range("J126").Value = 5161.125
$ws.Range("K126").Value = 11093.0001
$ws.Range("L126").Value = 15483.375
$ws.Range("M126").Value = -8623.000100000001
$ws.Range("N126").Value = -20423.375
$ws.Range("H132").Value = 1507.0952
$ws.Range("I132").Value = 1408.1875
$ws.Range("J132").Value = 1823.6
$ws.Range("K132").Value = 4224.5625
$ws.Range("L132").Value = 5470.799999999999
$ws.Range("M132").Value = -1694.5625
$ws.Range("N132").Value = -10530.8
$ws.Range("H134").Value = 1594.1578
$ws.Range("I134").Value = 1508.2941
$ws.Range("J134").Value = 2324
$ws.Range("K134").Value = 4524.8823
$ws.Range("L134").Value = 6972
$ws.Range("M134").Value = -1989.8823
$ws.Range("N134").Value = -12042
$ws.Range("H136").Value = 2630
$ws.Range("I136").Value = 2511.05
$ws.Range("K136").Value = 7533.150000000001
$ws.Range("M136").Value = -4983.150000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H21").Value = 387.33334
$ws.Range("I21").Value = 387.33334
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1162.00002
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -989.0000199999999
$ws.Range("N21").ClearContents()
$ws.Range("H26").Value = 96.666664
$ws.Range("I26").Value = 96
$ws.Range("J26").Value = 100
$ws.Range("K26").Value = 288
$ws.Range("L26").Value = 300
$ws.Range("M26").Value = 0
$ws.Range("N26").Value = -876
$ws.Range("H113").Value = 1739.4615
$ws.Range("I113").Value = 1590
$ws.Range("J113").Value = 1832.875
$ws.Range("K113").Value = 4770
$ws.Range("L113").Value = 5498.625
$ws.Range("M113").Value = -2600
$ws.Range("N113").Value = -9838.625
$ws.Range("H132").Value = 2499
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2499
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 22491
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -27551
$ws.Range("H137").Value = 2507.6924
$ws.Range("I137").Value = 933
$ws.Range("K137").Value = 2799
$ws.Range("M137").Value = 2301
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 328.42856
$ws.Range("I97").Value = 254.66667
$ws.Range("K97").Value = 254.66667
$ws.Range("M97").Value = 241.33333
$ws.Range("H102").Value = 2576.2083
$ws.Range("I102").Value = 2253.5217
$ws.Range("K102").Value = 2253.5217
$ws.Range("M102").Value = -631.5216999999998
$ws.Range("H122").Value = 28573782
$ws.Range("I122").Value = 40001988
$ws.Range("J122").Value = 3262
$ws.Range("K122").Value = 120005964
$ws.Range("L122").Value = 9786
$ws.Range("M122").Value = -120003514
$ws.Range("N122").Value = -14686
$ws.Range("H126").Value = 5247.8237
$ws.Range("I126").Value = 4621.273
$ws.Range("K126").Value = 13863.819
$ws.Range("M126").Value = -11393.819
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 925.8946999999999
$ws.Range("I16").Value = 1005.9375
$ws.Range("J16").Value = 499
$ws.Range("K16").Value = 1005.9375
$ws.Range("L16").Value = 499
$ws.Range("M16").Value = -835.9375
$ws.Range("N16").Value = -839
$ws.Range("H40").Value = 6614.1113
$ws.Range("I40").Value = 6627
$ws.Range("K40").Value = 6627
$ws.Range("M40").Value = -6491
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H93").Value = 5537.65
$ws.Range("I93").Value = 4995.933
$ws.Range("K93").Value = 4995.933
$ws.Range("M93").Value = -3747.933
$ws.Range("H132").Value = 5177.5713
$ws.Range("I132").Value = 4391.0713
$ws.Range("K132").Value = 13173.2139
$ws.Range("M132").Value = -10643.2139
$ws.Range("H136").Value = 5747
$ws.Range("I136").Value = 4644.8423
$ws.Range("J136").Value = 8364.625
$ws.Range("K136").Value = 13934.5269
$ws.Range("L136").Value = 25093.875
$ws.Range("M136").Value = -11384.5269
$ws.Range("N136").Value = -30193.875
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1111.4166
$ws.Range("I132").Value = 741.4815
$ws.Range("K132").Value = 2224.4445
$ws.Range("M132").Value = 305.5554999999999
